$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 112200000.0
$ws.Range("B4").Value = 28400000.0
$ws.Range("B9").Value = 24100000.0
$ws.Range("B10").Value = 504100000.0
$ws.Range("B11").Value = -197000000.0
$ws.Range("B16").Value = -5100000.0
$ws.Range("B17").Value = -202100000.0
$ws.Range("B18").Value = -2300000.0
$ws.Range("B19").Value = -122000000.0
$ws.Range("B20").Value = -48900000.0
$ws.Range("B22").Value = -173200000.0
$ws.Range("B23").Value = 19600000.0
$ws.Range("B24").Value = 148400000.0
$ws.Range("B25").Value = 615500000.0
$ws.Range("B26").Value = 763900000.0
$ws.Range("B27").Value = 34500000.0
$ws.Range("B28").Value = -48900000.0
$ws.Range("B29").Value = -59600000.0
$ws.Range("B31").Value = -122000000.0
